$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Acute pancreatitis - parenchymal edema" term (Pancreas row, row 14)
# to "Acute necrotizing pancreatitis " and add its YouTube link.
$ws.Range("B14").Value = "Acute necrotizing pancreatitis "

$ws.Hyperlinks.Add($ws.Range("D14"), "https://youtu.be/JvwODCASLYQ")
$ws.Range("D14").Value = "https://youtu.be/JvwODCASLYQ "
$ws.Range("D14").Style = "Collegamento ipertestuale"

# Move the active selection, as in the authored workbook.
$ws.Range("D19").Select()
